$d = $word.ActiveDocument

$newText = "Anda sedang berpartisipasi dalam kampanye global pengamatan dan pencatatan penampakan bintang paling redup untuk pengukuran tingkat polusi cahaya di suatu lokasi. Melalui pengamatan dan identifikasi  rasi bintang Hercules di langit malam dan membandingkannya dengan peta bintang, masyarakat di seluruh dunia dapat mengetahui dan mempelajari seberapa besar kontribusi cahaya di lingkungannya terhadap polusi cahaya. Kontribusi data anda pada basis data online akan membantu mendokumentasikan langit malam yang tampak di berbagai lokasi."

$rng = $d.Content
$found = $rng.Find.Execute("Anda sedang berpartisipasi*tampak di berbagai lokasi.", $true, $false, $true, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $rng.Delete()
    $rng.InsertAfter($newText)
}
